# "funcionalidades de la 61 a la 70 agregadas"
# Adds new test-report rows (IDs 61-90, worksheet rows 62-91) describing the
# "Creacion de Miembros" feature tests, merges the script/notes columns for
# the grouped-loop tests (rows 62-71), adds the explanatory rich-text note,
# and updates row 92 (ID 91) to reference Escenario-prueba91.js.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 62 (ID 61): first row of the merged block (B62:B71, C62:C71) ---
$ws.Range("B62").Value = "Escenario-prueba61-70.js"
$ws.Range("D62").Value = "Creacion de Miembros"
$ws.Range("E62").Value = "Creacion de miembro valido"
$ws.Range("G62").Value = "Martin Cantor"

# Rich text note in C62: plain intro paragraph + bold closing sentence
$introText = "Set a-priori generado con Mockaroo, descargado y adaptado`nEn la colección de datos (archivo de fixtures creacion_miembro.json) se tiene un registro por cada prueba a realizar y el resultado esperado de cada prueba, el script recore los datos y hace las 10 pruebas.`n"
$boldText = "Al poder hacer el loop, se facilitó la ejecución de escenarios diferentes en un solo script, por eso se agrupan y no existen los scripts del 61 al 70"
$ws.Range("C62").Value = $introText + $boldText
$ws.Range("C62").Characters($introText.Length + 1, $boldText.Length).Font.Bold = $true

# --- Rows 63-71 (IDs 62-70): only D/E/G differ, B/C come from the merge ---
$ws.Range("D63").Value = "Creacion de Miembros"
$ws.Range("E63").Value = "Creacion de miembro sin correo"
$ws.Range("G63").Value = "Martin Cantor"

$ws.Range("D64").Value = "Creacion de Miembros"
$ws.Range("E64").Value = "Creacion de miembro sin nombre"
$ws.Range("G64").Value = "Martin Cantor"

$ws.Range("D65").Value = "Creacion de Miembros"
$ws.Range("E65").Value = "Creacion de miembro sin Labels"
$ws.Range("G65").Value = "Martin Cantor"

$ws.Range("D66").Value = "Creacion de Miembros"
$ws.Range("E66").Value = "Creacion de miembro con correo  repetido"
$ws.Range("G66").Value = "Martin Cantor"

$ws.Range("D67").Value = "Creacion de Miembros"
$ws.Range("E67").Value = "Creacion de miembro con Note mayor a 500 caracteres"
$ws.Range("G67").Value = "Martin Cantor"

$ws.Range("D68").Value = "Creacion de Miembros"
$ws.Range("E68").Value = "Creacion de miembro con Note vacio"
$ws.Range("G68").Value = "Martin Cantor"

$ws.Range("D69").Value = "Creacion de Miembros"
$ws.Range("E69").Value = "Creacion de miembro y luego eliminarlo"
$ws.Range("G69").Value = "Martin Cantor"

$ws.Range("D70").Value = "Creacion de Miembros"
$ws.Range("E70").Value = "Creacion de miembro con correo invalido"
$ws.Range("G70").Value = "Martin Cantor"

$ws.Range("D71").Value = "Creacion de Miembros"
$ws.Range("E71").Value = "Creacion de miembro con Nombre corto"
$ws.Range("G71").Value = "Martin Cantor"
$ws.Rows.Item(71).RowHeight = 35

# Merge the grouped cells now that their content is set
$ws.Range("B62:B71").Merge()
$ws.Range("C62:C71").Merge()

# --- Rows 72-81 (IDs 71-80): individual scenario scripts, Mockaroo source ---
$mockarooSrc = "Set pseudo-aleatorio obtenido de Mockaroo por API, selección aleatoria del dato en el conjunto"
$fakerSrc = "tupla de datos generada dinámicamente durante la prueba usando faker"

$ws.Range("B72").Value = "Escenario-prueba71.js"
$ws.Range("C72").Value = $mockarooSrc
$ws.Range("D72").Value = "Creacion de Miembros"
$ws.Range("E72").Value = "Creacion de miembro con Note con saltos de linea"
$ws.Range("G72").Value = "Martin Cantor"

$ws.Range("B73").Value = "Escenario-prueba72.js"
$ws.Range("C73").Value = $mockarooSrc
$ws.Range("D73").Value = "Creacion de Miembros"
$ws.Range("E73").Value = "Creacion de miembro con correo sin arroba"
$ws.Range("G73").Value = "Martin Cantor"

$ws.Range("B74").Value = "Escenario-prueba73.js"
$ws.Range("C74").Value = $mockarooSrc
$ws.Range("D74").Value = "Creacion de Miembros"
$ws.Range("E74").Value = "Creacion de miembros con caracteres invalidos en el correo"
$ws.Range("G74").Value = "Martin Cantor"

$ws.Range("B75").Value = "Escenario-prueba74.js"
$ws.Range("C75").Value = $mockarooSrc
$ws.Range("D75").Value = "Creacion de Miembros"
$ws.Range("E75").Value = "Creacion de miembro con nombre muy largo"
$ws.Range("G75").Value = "Martin Cantor"

$ws.Range("B76").Value = "Escenario-prueba75.js"
$ws.Range("C76").Value = $mockarooSrc
$ws.Range("D76").Value = "Creacion de Miembros"
$ws.Range("E76").Value = "Creacion de miembro con correo con muchos caracteres"
$ws.Range("G76").Value = "Martin Cantor"

$ws.Range("B77").Value = "Escenario-prueba76.js"
$ws.Range("C77").Value = $mockarooSrc
$ws.Range("D77").Value = "Creacion de Miembros"
$ws.Range("E77").Value = "Creacion de miembro valido"
$ws.Range("G77").Value = "Martin Cantor"

$ws.Range("B78").Value = "Escenario-prueba77.js"
$ws.Range("C78").Value = $mockarooSrc
$ws.Range("D78").Value = "Creacion de Miembros"
$ws.Range("E78").Value = "Creacion de miembro sin correo"
$ws.Range("G78").Value = "Martin Cantor"

$ws.Range("B79").Value = "Escenario-prueba78.js"
$ws.Range("C79").Value = $mockarooSrc
$ws.Range("D79").Value = "Creacion de Miembros"
$ws.Range("E79").Value = "Creacion de miembro sin nombre"
$ws.Range("G79").Value = "Martin Cantor"

$ws.Range("B80").Value = "Escenario-prueba79.js"
$ws.Range("C80").Value = $mockarooSrc
$ws.Range("D80").Value = "Creacion de Miembros"
$ws.Range("E80").Value = "Creacion de miembro sin Labels"
$ws.Range("G80").Value = "Martin Cantor"

$ws.Range("B81").Value = "Escenario-prueba80.js"
$ws.Range("C81").Value = $mockarooSrc
$ws.Range("D81").Value = "Creacion de Miembros"
$ws.Range("E81").Value = "Creacion de miembro con correo  repetido"
$ws.Range("G81").Value = "Martin Cantor"

# --- Rows 82-91 (IDs 81-90): faker-sourced tuples ---
$ws.Range("B82").Value = "Escenario-prueba80.js"
$ws.Range("C82").Value = $fakerSrc
$ws.Range("D82").Value = "Creacion de Miembros"
$ws.Range("E82").Value = "Creacion de miembro con Note mayor a 500 caracteres"
$ws.Range("G82").Value = "Martin Cantor"

$ws.Range("B83").Value = "Escenario-prueba80.js"
$ws.Range("C83").Value = $fakerSrc
$ws.Range("D83").Value = "Creacion de Miembros"
$ws.Range("E83").Value = "Creacion de miembro con Note vacio"
$ws.Range("G83").Value = "Martin Cantor"

$ws.Range("B84").Value = "Escenario-prueba80.js"
$ws.Range("C84").Value = $fakerSrc
$ws.Range("D84").Value = "Creacion de Miembros"
$ws.Range("E84").Value = "Creacion de miembro y luego eliminarlo"
$ws.Range("G84").Value = "Martin Cantor"

$ws.Range("B85").Value = "Escenario-prueba80.js"
$ws.Range("C85").Value = $fakerSrc
$ws.Range("D85").Value = "Creacion de Miembros"
$ws.Range("E85").Value = "Creacion de miembro con correo invalido"
$ws.Range("G85").Value = "Martin Cantor"

$ws.Range("B86").Value = "Escenario-prueba80.js"
$ws.Range("C86").Value = $fakerSrc
$ws.Range("D86").Value = "Creacion de Miembros"
$ws.Range("E86").Value = "Creacion de miembro con Nombre corto"
$ws.Range("G86").Value = "Martin Cantor"

$ws.Range("B87").Value = "Escenario-prueba80.js"
$ws.Range("C87").Value = $fakerSrc
$ws.Range("D87").Value = "Creacion de Miembros"
$ws.Range("E87").Value = "Creacion de miembro con Note con saltos de linea"
$ws.Range("G87").Value = "Martin Cantor"

$ws.Range("B88").Value = "Escenario-prueba80.js"
$ws.Range("C88").Value = $fakerSrc
$ws.Range("D88").Value = "Creacion de Miembros"
$ws.Range("E88").Value = "Creacion de miembro con correo sin arroba"
$ws.Range("G88").Value = "Martin Cantor"

$ws.Range("B89").Value = "Escenario-prueba80.js"
$ws.Range("C89").Value = $fakerSrc
$ws.Range("D89").Value = "Creacion de Miembros"
$ws.Range("E89").Value = "Creacion de miembros con caracteres invalidos en el correo"
$ws.Range("G89").Value = "Martin Cantor"

$ws.Range("B90").Value = "Escenario-prueba80.js"
$ws.Range("C90").Value = $fakerSrc
$ws.Range("D90").Value = "Creacion de Miembros"
$ws.Range("E90").Value = "Creacion de miembro con nombre muy largo"
$ws.Range("G90").Value = "Martin Cantor"

$ws.Range("B91").Value = "Escenario-prueba90.js"
$ws.Range("C91").Value = $fakerSrc
$ws.Range("D91").Value = "Creacion de Miembros"
$ws.Range("E91").Value = "Creacion de miembro con correo con muchos caracteres"
$ws.Range("G91").Value = "Martin Cantor"

# --- Row 92 (ID 91): now references Escenario-prueba91.js instead of 90.js ---
$ws.Range("B92").Value = "Escenario-prueba91.js"
